$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings (e.g. leading/trailing zeros, thousand-dot
# separators) that must stay text. Force the text number format before writing
# values that Excel would otherwise auto-convert to a numeric value.
$ws.Range('D2').Value = '36.457.58'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.939.71'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.82'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.610'
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.96'
$ws.Range('E8').Value = '  -3.65%  '
$ws.Range('E9').Value = '  -3.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0852'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '2.224.64'
$ws.Range('E12').Value = '  -1.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.32'
$ws.Range('E13').Value = '  -6.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.807'
$ws.Range('E14').Value = '  -6.01%  '
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('E16').Value = '  -5.92%  '
$ws.Range('D17').Value = '1.940.08'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '36.408.79'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0864'
$ws.Range('E19').Value = '  -3.55%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.13'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '227.30'
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.98'
$ws.Range('E22').Value = '  -5.59%  '
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.33'
$ws.Range('E24').Value = '  -6.83%  '
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.15'
$ws.Range('E26').Value = '  -7.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.00'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.21'
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('E30').Value = '  -1.57%  '
$ws.Range('E31').Value = '  -7.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.55'
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('E33').Value = '  -6.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.16'
$ws.Range('E34').Value = '  -6.78%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.06'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.79'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.18'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.06'
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0991'
$ws.Range('E40').Value = '  +2.86%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('E43').Value = '  -5.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.58'
$ws.Range('E44').Value = '  -3.65%  '
$ws.Range('D45').Value = '1.340.44'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('E46').Value = '  -6.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.18'
$ws.Range('E47').Value = '  -5.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.10'
$ws.Range('E48').Value = '  -4.60%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = '2.116.30'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.20'
$ws.Range('E51').Value = '  -4.87%  '
